$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header block ---
$ws.Range("C2").Value = "Hartmut"

# B3 holds a long, purely-numeric card number that must stay text (inlineStr in
# the source). Assigning a numeric-looking string via .Value auto-converts it
# to a number, so force Text format first, then restore the original style
# (s=8) by pasting formats from a cell that already carries that style.
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C2").Copy()
$ws.Range("B3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("C3").Value = "Mohaupt"
$ws.Range("D5").Value = "KONTOSTAND AM 29.09.2024"

# --- Row 6 ---
$ws.Range("B6").Value = "01.10."
$ws.Range("C6").Value = "02.10."
$ws.Range("D6").Value = "PAYPAL SKTISH"
$ws.Range("E6").Value = "31,80-"

# --- Row 7 ---
$ws.Range("B7").Value = "05.10."
$ws.Range("C7").Value = "06.10."
$ws.Range("D7").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 96955554"
$ws.Range("E7").Value = "83,61-"

# --- Row 8 ---
$ws.Range("B8").Value = "09.10."
$ws.Range("C8").Value = "10.10."
$ws.Range("D8").Value = "EBAY MKTPLC EU YRCNJN"
$ws.Range("E8").Value = "180,35-"

# --- Row 9 ---
$ws.Range("B9").Value = "12.10."
$ws.Range("C9").Value = "13.10."
$ws.Range("D9").Value = "BEITRAG Allianz SE K-93811780"
$ws.Range("E9").Value = "56,67-"

# --- Row 10 (previously blank, now populated) ---
$ws.Range("B10").Value = "15.10."
$ws.Range("C10").Value = "16.10."
$ws.Range("D10").Value = "RECHNUNG VODAFONE GMBH 90275252"
$ws.Range("E10").Value = "39,90-"
$ws.Range("E10").HorizontalAlignment = -4152  # xlHAlignRight
$ws.Range("E10").VerticalAlignment = -4107    # xlVAlignBottom
$ws.Range("E10").WrapText = $false

# --- Row 11 (previously blank, now populated) ---
$ws.Range("B11").Value = "16.10."
$ws.Range("C11").Value = "17.10."
$ws.Range("D11").Value = "PAYPAL JCRDYW"
$ws.Range("E11").Value = "7,04-"
$ws.Range("E11").HorizontalAlignment = -4152  # xlHAlignRight
$ws.Range("E11").VerticalAlignment = -4107    # xlVAlignBottom
$ws.Range("E11").WrapText = $false

# --- Row 12 ---
$ws.Range("D12").Value = "KONTOSTAND AM 20.10.2024"
$ws.Range("E12").Value = "399,37-"

# --- Row 13 ---
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 26.10.2024"
